# Auto update stock data
# Updates the "as of" date and EBITDA figures for the latest snapshot rows.
# Values are written with a leading apostrophe so Excel stores them as
# literal text (matching the original inlineStr/text cell type) instead of
# auto-converting them into date serials / numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Date (column A) and new EBITDA (column B, optional)
$updates = @(
    @{ Row = 2;  Date = "2025/11/27"; Ebitda = "5.09" }
    @{ Row = 8;  Date = "2025/11/27"; Ebitda = "7.72" }
    @{ Row = 14; Date = "2025/11/27"; Ebitda = "2.79" }
    @{ Row = 20; Date = "2025/11/27"; Ebitda = "12.42" }
    @{ Row = 26; Date = "2025/11/27"; Ebitda = "9.96" }
    @{ Row = 32; Date = "2025/11/27"; Ebitda = "26.04" }
    @{ Row = 38; Date = "2025/11/27"; Ebitda = $null }
    @{ Row = 44; Date = "2025/11/27"; Ebitda = "10.90" }
    @{ Row = 50; Date = "2025/11/27"; Ebitda = "11.61" }
    @{ Row = 56; Date = "2025/11/27"; Ebitda = "34.78" }
    @{ Row = 62; Date = "2025/11/27"; Ebitda = "11.28" }
    @{ Row = 68; Date = "2025/11/27"; Ebitda = "12.23" }
    @{ Row = 74; Date = "2025/11/27"; Ebitda = "15.74" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = "'" + $u.Date
    if ($u.Ebitda) {
        $ws.Cells.Item($u.Row, 2).Value = "'" + $u.Ebitda
    }
}
